$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "environ" column (E) with 1 for each data row (2-7)
$ws.Range("E2:E7").Value = 1

# Move the active selection to E8, matching the saved view state
$ws.Range("E8").Select()
